$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = "Abdullah Al Kalbani"
$summary.Range("B4").Value = 1896.28
$summary.Range("B6").Value = 3523
$summary.Range("B7").Value = 24686
$summary.Range("B8").Value = -21163
$summary.Range("B9").Value = 0.14

# --- Sheet "Assets" ---
# Remove the "Vehicles / Mid-range Car" row (row 2), shifting "Liquid Assets" row up
$assets = $wb.Worksheets.Item("Assets")
$assets.Rows.Item(2).Delete()
# Now row 2 is "Liquid Assets / Savings Account", row 3 is "TOTAL ASSETS"
$assets.Range("C2").Value = 3523
$assets.Range("C3").Value = 3523

# --- Sheet "Liabilities" ---
# Remove "Auto Loans" (row 2) and "Personal Loans" (row 3), shifting "Credit Cards" up
$liabilities = $wb.Worksheets.Item("Liabilities")
$liabilities.Rows.Item(2).Delete()
$liabilities.Rows.Item(2).Delete()
# Now row 2 is "Credit Cards / Credit Card Balance", row 3 is "TOTAL LIABILITIES"
$liabilities.Range("C2").Value = 24686
$liabilities.Range("D2").Value = 1234
$liabilities.Range("E2").Value = 1
$liabilities.Range("C3").Value = 24686
